$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ECs" row (original row 2) and the "Inflammatory-Mac" row
# (original row 4, which becomes row 3 after the first deletion), leaving
# only the "FAPs" and "MuSCs" sending-cluster rows, now updated with new
# TPM-derived values.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(3).Delete()

# Row 2: FAPs -> Gm13306 -> Ackr2 -> FAPs
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Gm13306"
$ws.Cells.Item(2, 3).Value = "Ackr2"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.1234563333333333
$ws.Cells.Item(2, 8).Value = 0.370369
$ws.Cells.Item(2, 9).Value = 0.4718467134221305
$ws.Cells.Item(2, 10).Value = 0.4718467134221305
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.062395
$ws.Cells.Item(2, 14).Value = 0.187185
$ws.Cells.Item(2, 15).Value = 1
$ws.Cells.Item(2, 16).Value = 1
$ws.Cells.Item(2, 17).Value = 0.007703057918333333
$ws.Cells.Item(2, 18).Value = 0.069327521265
$ws.Cells.Item(2, 19).Value = 0.4718467134221305
$ws.Cells.Item(2, 20).Value = 0.4718467134221305

# Row 3: MuSCs -> Gm13306 -> Ackr2 -> FAPs
$ws.Cells.Item(3, 1).Value = "MuSCs"
$ws.Cells.Item(3, 2).Value = "Gm13306"
$ws.Cells.Item(3, 3).Value = "Ackr2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.1381886666666667
$ws.Cells.Item(3, 8).Value = 0.414566
$ws.Cells.Item(3, 9).Value = 0.5281532865778694
$ws.Cells.Item(3, 10).Value = 0.5281532865778695
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.062395
$ws.Cells.Item(3, 14).Value = 0.187185
$ws.Cells.Item(3, 15).Value = 1
$ws.Cells.Item(3, 16).Value = 1
$ws.Cells.Item(3, 17).Value = 0.008622281856666665
$ws.Cells.Item(3, 18).Value = 0.07760053671
$ws.Cells.Item(3, 19).Value = 0.5281532865778694
$ws.Cells.Item(3, 20).Value = 0.5281532865778695

Write-Output "Done. UsedRange rows: $($ws.UsedRange.Rows.Count) cols: $($ws.UsedRange.Columns.Count)"
